$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Rename the sheet ---
$ws.Name = "Dic_Einheit"

# --- Update column widths ---
$ws.Columns.Item(1).ColumnWidth = 22.8515625
$ws.Columns.Item(2).ColumnWidth = 67.52734375
$ws.Columns.Item(3).ColumnWidth = 82.6171875

# --- Extend formatting (style of row 16, which mirrors the header-data style) down to the new rows 17-28 ---
$ws.Range("A16:C16").Copy()
$ws.Range("A17:C28").PasteSpecial(-4122)

# --- Write header + data rows (plain text values) ---
$ws.Range("A1").Value = "ENr"
$ws.Range("B1").Value = "Einheit De"
$ws.Range("C1").Value = "Einheit En"
$ws.Range("A2").Value = "E_10H3"
$ws.Range("B2").Value = "1 000"
$ws.Range("A3").Value = "E_10H7"
$ws.Range("B3").Value = "Millionen"
$ws.Range("C3").Value = "Millions"
$ws.Range("A4").Value = "E_BNEUR"
$ws.Range("B4").Value = "Milliarden EUR"
$ws.Range("C4").Value = "Billion EUR"
$ws.Range("A5").Value = "E_EWPKM2"
$ws.Range("B5").Value = "Einwohner/-innen pro m² Siedlungs- und Verkehrsfläche"
$ws.Range("C5").Value = "Inhabitants per m² settlement and transport area"
$ws.Range("A6").Value = "E_GINI"
$ws.Range("B6").Value = "Gini-Koeffizient"
$ws.Range("C6").Value = "Gini coefficient"
$ws.Range("A7").Value = "E_HAPD"
$ws.Range("B7").Value = "Hektar pro Tag"
$ws.Range("C7").Value = "Hectre per day"
$ws.Range("A8").Value = "E_IDX_1990"
$ws.Range("B8").Value = "1990 = 100"
$ws.Range("C8").Value = "1990 = 100"
$ws.Range("A9").Value = "E_IDX_2030"
$ws.Range("B9").Value = "2030 = 100"
$ws.Range("C9").Value = "2030 = 100"
$ws.Range("A10").Value = "E_IDX2000"
$ws.Range("B10").Value = "2000 = 100"
$ws.Range("C10").Value = "2000 = 100"
$ws.Range("A11").Value = "E_IDX2005"
$ws.Range("B11").Value = "2005 = 100"
$ws.Range("C11").Value = "2005 = 100"
$ws.Range("A12").Value = "E_IDX2008"
$ws.Range("B12").Value = "2008 = 100"
$ws.Range("C12").Value = "2008 = 100"
$ws.Range("A13").Value = "E_IDX2010"
$ws.Range("B13").Value = "2010 = 100"
$ws.Range("C13").Value = "2010 = 100"
$ws.Range("A14").Value = "E_IDX2015"
$ws.Range("B14").Value = "2015 = 100"
$ws.Range("C14").Value = "2015 = 100"
$ws.Range("A15").Value = "E_KGPHA"
$ws.Range("B15").Value = "Kilogramm pro Hektar"
$ws.Range("C15").Value = "Kilogram per hectare"
$ws.Range("A16").Value = "E_MILIGPL"
$ws.Range("B16").Value = "Milligramm pro Liter"
$ws.Range("C16").Value = "Miligrams per litre"
$ws.Range("A17").Value = "E_MIN"
$ws.Range("B17").Value = "Minuten"
$ws.Range("C17").Value = "Minutes"
$ws.Range("A18").Value = "E_MNEUR"
$ws.Range("B18").Value = "Millionen EUR"
$ws.Range("C18").Value = "Million EUR"
$ws.Range("A19").Value = "E_MNEW"
$ws.Range("B19").Value = "Millionen Einwohner/-innen"
$ws.Range("C19").Value = "Million inhabitants"
$ws.Range("A20").Value = "E_NN"
$ws.Range("A21").Value = "E_NUM"
$ws.Range("B21").Value = "Anzahl"
$ws.Range("C21").Value = "Number"
$ws.Range("A22").Value = "E_P10H6EWN"
$ws.Range("B22").Value = "Je 100 000 Einwohner/-innen"
$ws.Range("C22").Value = "Per 100,000 inhabitants"
$ws.Range("A23").Value = "E_P10H6EWNU70"
$ws.Range("B23").Value = "Je 100 000 Einwohner/-innen unter 70 Jahren (ohne unter 1-Jährige)"
$ws.Range("C23").Value = "Per 100,000 inhabitants under 70 years (excluding under 1 year olds)"
$ws.Range("A24").Value = "E_PRZNT"
$ws.Range("B24").Value = "Prozent"
$ws.Range("C24").Value = "Percentage"
$ws.Range("A25").Value = "E_PRZNTPKT"
$ws.Range("B25").Value = "Prozentpunkte"
$ws.Range("C25").Value = "Percentage points"
$ws.Range("A26").Value = "E_QMPA"
$ws.Range("B26").Value = "m² pro Jahr"
$ws.Range("C26").Value = "m² per year"
$ws.Range("A27").Value = "E_QMPINHABA"
$ws.Range("B27").Value = "m² pro Einwohner/-in und Jahr"
$ws.Range("C27").Value = "m² per inhabitant and year"
$ws.Range("A28").Value = "E_TEUR"
$ws.Range("B28").Value = "1 000 EUR"
$ws.Range("C28").Value = "1.000 EUR"

# --- Special-case cells whose literal text would otherwise be auto-coerced ---
# C2 = "1.000" looks numeric to the auto-type-detection; force literal text via
# a leading apostrophe, then restore the plain (non quote-prefixed) style of the column.
$ws.Range("C2").Value = "'1.000"
$ws.Range("C16").Copy()
$ws.Range("C2").PasteSpecial(-4122)

# B20/C20 (E_NN row) are empty-string cells, not blank/unset cells -
# an apostrophe forces a genuine (empty) text cell instead of clearing it,
# then restore the rows plain style.
$ws.Range("B20").Value = "'"
$ws.Range("C20").Value = "'"
$ws.Range("C16").Copy()
$ws.Range("B20").PasteSpecial(-4122)
$ws.Range("C16").Copy()
$ws.Range("C20").PasteSpecial(-4122)
